$d = $word.ActiveDocument

$d.Content.Find.Execute("857÷6=142, 5", $true, $false, $false, $false, $false, $true, 1, $false, "614÷8=76, 6", 2) | Out-Null
$d.Content.Find.Execute("328÷3=109, 1", $true, $false, $false, $false, $false, $true, 1, $false, "990÷2=495, 0", 2) | Out-Null
$d.Content.Find.Execute("605÷9=67, 2", $true, $false, $false, $false, $false, $true, 1, $false, "404÷2=202, 0", 2) | Out-Null
$d.Content.Find.Execute("995÷9=110, 5", $true, $false, $false, $false, $false, $true, 1, $false, "201÷6=33, 3", 2) | Out-Null
$d.Content.Find.Execute("880÷2=440, 0", $true, $false, $false, $false, $false, $true, 1, $false, "709÷6=118, 1", 2) | Out-Null
$d.Content.Find.Execute("489÷5=97, 4", $true, $false, $false, $false, $false, $true, 1, $false, "453÷2=226, 1", 2) | Out-Null
$d.Content.Find.Execute("892÷4=223, 0", $true, $false, $false, $false, $false, $true, 1, $false, "431÷4=107, 3", 2) | Out-Null
$d.Content.Find.Execute("360÷2=180, 0", $true, $false, $false, $false, $false, $true, 1, $false, "438÷8=54, 6", 2) | Out-Null
$d.Content.Find.Execute("552÷9=61, 3", $true, $false, $false, $false, $false, $true, 1, $false, "443÷3=147, 2", 2) | Out-Null
$d.Content.Find.Execute("781÷9=86, 7", $true, $false, $false, $false, $false, $true, 1, $false, "124÷9=13, 7", 2) | Out-Null
$d.Content.Find.Execute("398÷6=66, 2", $true, $false, $false, $false, $false, $true, 1, $false, "811÷3=270, 1", 2) | Out-Null
$d.Content.Find.Execute("358÷2=179, 0", $true, $false, $false, $false, $false, $true, 1, $false, "968÷2=484, 0", 2) | Out-Null
$d.Content.Find.Execute("963÷7=137, 4", $true, $false, $false, $false, $false, $true, 1, $false, "894÷3=298, 0", 2) | Out-Null
$d.Content.Find.Execute("649÷7=92, 5", $true, $false, $false, $false, $false, $true, 1, $false, "282÷8=35, 2", 2) | Out-Null
$d.Content.Find.Execute("484÷2=242, 0", $true, $false, $false, $false, $false, $true, 1, $false, "584÷2=292, 0", 2) | Out-Null
$d.Content.Find.Execute("916÷7=130, 6", $true, $false, $false, $false, $false, $true, 1, $false, "383÷5=76, 3", 2) | Out-Null
$d.Content.Find.Execute("377÷8=47, 1", $true, $false, $false, $false, $false, $true, 1, $false, "742÷2=371, 0", 2) | Out-Null
$d.Content.Find.Execute("785÷7=112, 1", $true, $false, $false, $false, $false, $true, 1, $false, "411÷7=58, 5", 2) | Out-Null
$d.Content.Find.Execute("547÷6=91, 1", $true, $false, $false, $false, $false, $true, 1, $false, "534÷9=59, 3", 2) | Out-Null
$d.Content.Find.Execute("783÷2=391, 1", $true, $false, $false, $false, $false, $true, 1, $false, "807÷6=134, 3", 2) | Out-Null
$d.Content.Find.Execute("883÷9=98, 1", $true, $false, $false, $false, $false, $true, 1, $false, "444÷4=111, 0", 2) | Out-Null
$d.Content.Find.Execute("117÷3=39, 0", $true, $false, $false, $false, $false, $true, 1, $false, "899÷8=112, 3", 2) | Out-Null
$d.Content.Find.Execute("701÷3=233, 2", $true, $false, $false, $false, $false, $true, 1, $false, "922÷4=230, 2", 2) | Out-Null
$d.Content.Find.Execute("546÷8=68, 2", $true, $false, $false, $false, $false, $true, 1, $false, "776÷2=388, 0", 2) | Out-Null
$d.Content.Find.Execute("637÷5=127, 2", $true, $false, $false, $false, $false, $true, 1, $false, "629÷9=69, 8", 2) | Out-Null
